# Add the new "2022-Q4" quarterly sheet right after "总计" (i.e. right before
# the existing "2022-Q3" sheet), and update the "总计" summary sheet with a
# row for it.

$wb = $excel.ActiveWorkbook

# --- 1. Insert a new worksheet right before the current "2022-Q3" sheet ---
$sheetQ3Before = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($sheetQ3Before, $null)
$newSheet.Name = "2022-Q4"

# Re-fetch the "2022-Q3" sheet by name now that the insert/rename is done --
# the handle obtained before the Add() call tracks the new sheet's slot, not
# the original sheet, once a sheet has been inserted in front of it.
$sheetQ3 = $wb.Worksheets.Item("2022-Q3")

# Pull header-row + index-cell formatting from an existing quarter sheet so
# the new sheet matches the look of its siblings.
$sheetQ3.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$sheetQ3.Range("A2").Copy($newSheet.Range("A2"))

# Fund holdings for 2022-Q4 (single fund on record for this quarter).
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Cells.Item(2, 2).Value = "320022"
$newSheet.Cells.Item(2, 3).Value = "诺安研究精选股票"
$newSheet.Cells.Item(2, 4).Value = "6.42"
$newSheet.Cells.Item(2, 5).Value = "92.87"
$newSheet.Cells.Item(2, 6).Value = "2.15"
$newSheet.Cells.Item(2, 7).Value = "0.1380"
$newSheet.Cells.Item(2, 8).Value = 7

# --- 2. Update the "总计" summary sheet: add the 2022-Q4 row at the top of
#        the data (row 2), shifting every other quarter down by one. ---
$total = $wb.Worksheets.Item("总计")

$quarters = @(
    @("2022-Q4", 1, 0.14),
    @("2022-Q3", 2, 0.24),
    @("2022-Q2", 9, 2.46),
    @("2022-Q1", 36, 10.84),
    @("2021-Q4", 43, 27.39),
    @("2021-Q3", 30, 14.46),
    @("2021-Q2", 28, 8.85),
    @("2021-Q1", 23, 11.09),
    @("2020-Q4", 6, 1.51)
)

# Row 10 is brand new (the sheet used to stop at row 9) -- seed its "A"
# index cell from row 9's so it picks up the same index-column style.
$total.Range("A9").Copy($total.Range("A10"))

for ($i = 0; $i -lt $quarters.Length; $i++) {
    $row = 2 + $i
    $entry = $quarters[$i]
    $total.Cells.Item($row, 1).Value = $i
    $total.Cells.Item($row, 2).Value = $entry[0]
    $total.Cells.Item($row, 3).Value = $entry[1]
    $total.Cells.Item($row, 4).Value = $entry[2]
}
